$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "ip_address_list" - re-sorted / re-edited project list
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ip_address_list")

# Clear the previously existing range first (old data spanned A1:E15)
$ws1.Range("A1:E15").Clear()

# Row 1
$ws1.Cells.Item(1,1).Value = "kkkk"
$ws1.Cells.Item(1,2).Value = "10.96.205.17"
$ws1.Cells.Item(1,3).Value = "255.255.255.0"
$ws1.Cells.Item(1,4).Value = "PC:`t10.96.205.175NAS:`t10.96.205.166k`nFH:`t10.96.205.154`n`t10.96.20aa"
$ws1.Cells.Item(1,5).Value = 0

# Row 2
$ws1.Cells.Item(2,1).Value = "jjjjb"
$ws1.Cells.Item(2,2).Value = "192.168.000.000"
$ws1.Cells.Item(2,3).Value = "255.255.255.0"
$ws1.Cells.Item(2,4).Value = "jh"
$ws1.Cells.Item(2,5).Value = 0

# Row 3
$ws1.Cells.Item(3,1).Value = "hhhd"
$ws1.Cells.Item(3,2).Value = "192.168.000.000"
$ws1.Cells.Item(3,3).Value = "255.255.255.0"
$ws1.Cells.Item(3,5).Value = $false

# Row 4
$ws1.Cells.Item(4,1).Value = "axggg"
$ws1.Cells.Item(4,2).Value = "192.168.000.000xg"
$ws1.Cells.Item(4,3).Value = "255.255.255.0"
$ws1.Cells.Item(4,4).Value = "axg"
$ws1.Cells.Item(4,5).Value = $true

# Row 5
$ws1.Cells.Item(5,1).Value = "hhgggg"
$ws1.Cells.Item(5,2).Value = "192.168.000.000h"
$ws1.Cells.Item(5,3).Value = "255.255.255.0"
$ws1.Cells.Item(5,4).Value = "hhh"
$ws1.Cells.Item(5,5).Value = $true

# Row 6
$ws1.Cells.Item(6,1).Value = "jjs"
$ws1.Cells.Item(6,2).Value = "192.168.000.000j"
$ws1.Cells.Item(6,3).Value = "255.255.255.0"
$ws1.Cells.Item(6,4).Value = "ssjdg"
$ws1.Cells.Item(6,5).Value = $true

# Row 7
$ws1.Cells.Item(7,1).Value = "sega"
$ws1.Cells.Item(7,2).Value = "192.168.000.000"
$ws1.Cells.Item(7,3).Value = "255.255.255.0"
$ws1.Cells.Item(7,4).Value = "g"
$ws1.Cells.Item(7,5).Value = $false

# Row 8
$ws1.Cells.Item(8,1).Value = "aaujxa"
$ws1.Cells.Item(8,2).Value = "192.168.000.000"
$ws1.Cells.Item(8,3).Value = "255.255.255.0"
$ws1.Cells.Item(8,5).Value = 1

# Row 9
$ws1.Cells.Item(9,1).Value = "ggzagga"
$ws1.Cells.Item(9,2).Value = "192.168.000.0g"
$ws1.Cells.Item(9,3).Value = "255.255.255.0"
$ws1.Cells.Item(9,4).Value = "gg"
$ws1.Cells.Item(9,5).Value = $true

# Row 10
$ws1.Cells.Item(10,1).Value = "haxs"
$ws1.Cells.Item(10,2).Value = "192.168.000.000"
$ws1.Cells.Item(10,3).Value = "255.255.255.0"
$ws1.Cells.Item(10,4).Value = "x"
$ws1.Cells.Item(10,5).Value = 1

# ---------------------------------------------------------------
# Sheet "ip_adress_fav_list" - re-sorted / re-edited favourites
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")

$ws2.Range("A1:E7").Clear()

# Row 1
$ws2.Cells.Item(1,1).Value = "axggg"
$ws2.Cells.Item(1,2).Value = "192.168.000.000xg"
$ws2.Cells.Item(1,3).Value = "255.255.255.0"
$ws2.Cells.Item(1,4).Value = "axg"
$ws2.Cells.Item(1,5).Value = $true

# Row 2
$ws2.Cells.Item(2,1).Value = "hhgggg"
$ws2.Cells.Item(2,2).Value = "192.168.000.000h"
$ws2.Cells.Item(2,3).Value = "255.255.255.0"
$ws2.Cells.Item(2,4).Value = "hhh"
$ws2.Cells.Item(2,5).Value = $true

# Row 3
$ws2.Cells.Item(3,1).Value = "jjs"
$ws2.Cells.Item(3,2).Value = "192.168.000.000j"
$ws2.Cells.Item(3,3).Value = "255.255.255.0"
$ws2.Cells.Item(3,4).Value = "ssjdg"
$ws2.Cells.Item(3,5).Value = $true

# Row 4
$ws2.Cells.Item(4,1).Value = "aaujxa"
$ws2.Cells.Item(4,2).Value = "192.168.000.000"
$ws2.Cells.Item(4,3).Value = "255.255.255.0"
$ws2.Cells.Item(4,5).Value = 1

# Row 5
$ws2.Cells.Item(5,1).Value = "ggzagga"
$ws2.Cells.Item(5,2).Value = "192.168.000.0g"
$ws2.Cells.Item(5,3).Value = "255.255.255.0"
$ws2.Cells.Item(5,4).Value = "gg"
$ws2.Cells.Item(5,5).Value = $true

# Row 6
$ws2.Cells.Item(6,1).Value = "haxs"
$ws2.Cells.Item(6,2).Value = "192.168.000.000"
$ws2.Cells.Item(6,3).Value = "255.255.255.0"
$ws2.Cells.Item(6,4).Value = "x"
$ws2.Cells.Item(6,5).Value = 1

# ---------------------------------------------------------------
# Sheet "Settings" - update a couple of flags, add new row
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Settings")
$ws4.Cells.Item(3,2).Value = 1
$ws4.Cells.Item(8,2).Value = 0
$ws4.Cells.Item(9,2).Value = 0
